# Atualização de bases das ligas, do dia: 18-04-2024 às 00:36
# Uruguay Primera Divisão - corrige a ordem de 2 pares de partidas (114/115, 117/118),
# atualiza os resultados/odds das partidas 181/182 e adiciona as novas partidas 183/184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($row, $data)
    foreach ($col in $data.Keys) {
        $ws.Range($col + $row).Value = $data[$col]
    }
}

# --- Row 114: refresh match data (id/style of column A stay put) ---
$d114 = [ordered]@{
    "B" = 7559469
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Clausura"
    "E" = 45266.70833333334
    "F" = "Montevideo Wanderers"
    "G" = "Penarol"
    "H" = 0
    "I" = 0
    "J" = "D"
    "K" = 4.75
    "L" = 3.4
    "M" = 1.7
    "N" = 2.7
    "O" = 3.2
    "P" = 2.45
    "Q" = 0
    "R" = 2.05
    "S" = 1.8
    "T" = 2.5
    "U" = 1.975
    "V" = 1.875
    "W" = -1
    "X" = 2.2
    "Y" = -1
    "Z" = 0
    "AA" = 0
    "AB" = -1
    "AC" = 0.875
}
Set-RowData 114 $d114

# --- Row 115: refresh match data (id/style of column A stay put) ---
$d115 = [ordered]@{
    "B" = 7559468
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Clausura"
    "E" = 45266.70833333334
    "F" = "Liverpool Montevideo"
    "G" = "CA River Plate"
    "H" = 2
    "I" = 1
    "J" = "H"
    "K" = 1.7
    "L" = 3
    "M" = 5.75
    "N" = 1.833
    "O" = 3.2
    "P" = 4.5
    "Q" = -0.5
    "R" = 1.925
    "S" = 1.925
    "T" = 2.25
    "U" = 2.025
    "V" = 1.825
    "W" = 0.833
    "X" = -1
    "Y" = -1
    "Z" = 0.925
    "AA" = -1
    "AB" = 1.025
    "AC" = -1
}
Set-RowData 115 $d115

# --- Row 117: refresh match data (id/style of column A stay put) ---
$d117 = [ordered]@{
    "B" = 7013885
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Clausura"
    "E" = 45267.70833333334
    "F" = "La Luz"
    "G" = "Atletico Fenix Montevideo"
    "H" = 0
    "I" = 2
    "J" = "A"
    "K" = 3
    "L" = 3
    "M" = 2.4
    "N" = 2.9
    "O" = 2.75
    "P" = 2.6
    "Q" = 0
    "R" = 2.025
    "S" = 1.825
    "T" = 2
    "U" = 2.025
    "V" = 1.825
    "W" = -1
    "X" = -1
    "Y" = 1.6
    "Z" = -1
    "AA" = 0.825
    "AB" = 0
    "AC" = 0
}
Set-RowData 117 $d117

# --- Row 118: refresh match data (id/style of column A stay put) ---
$d118 = [ordered]@{
    "B" = 7013702
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Clausura"
    "E" = 45267.70833333334
    "F" = "Defensor Sporting"
    "G" = "Danubio"
    "H" = 0
    "I" = 2
    "J" = "A"
    "K" = 1.8
    "L" = 3.6
    "M" = 4.2
    "N" = 1.8
    "O" = 3.6
    "P" = 4.2
    "Q" = -0.75
    "R" = 2.05
    "S" = 1.8
    "T" = 2.25
    "U" = 1.85
    "V" = 2
    "W" = -1
    "X" = -1
    "Y" = 3.2
    "Z" = -1
    "AA" = 0.8
    "AB" = -0.5
    "AC" = 0.5
}
Set-RowData 118 $d118

# --- Row 181: refresh match data (id/style of column A stay put) ---
$d181 = [ordered]@{
    "B" = 8051187
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Apertura"
    "E" = 45396.625
    "F" = "Defensor Sporting"
    "G" = "CA River Plate"
    "H" = 2
    "I" = 1
    "J" = "H"
    "K" = 1.727
    "L" = 3.5
    "M" = 5
    "N" = 1.6
    "O" = 3.6
    "P" = 6
    "Q" = -0.75
    "R" = 1.8
    "S" = 2.05
    "T" = 2.25
    "U" = 1.875
    "V" = 1.975
    "W" = 0.6000000000000001
    "X" = -1
    "Y" = -1
    "Z" = 0.4
    "AA" = -0.5
    "AB" = 0.875
    "AC" = -1
}
Set-RowData 181 $d181

# --- Row 182: refresh match data (id/style of column A stay put) ---
$d182 = [ordered]@{
    "B" = 8050911
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Apertura"
    "E" = 45396.75
    "F" = "Penarol"
    "G" = "Danubio"
    "H" = 2
    "I" = 0
    "J" = "H"
    "K" = 1.666
    "L" = 3.5
    "M" = 5.5
    "N" = 1.6
    "O" = 3.6
    "P" = 6.5
    "Q" = -0.75
    "R" = 1.8
    "S" = 2.05
    "T" = 2.25
    "U" = 2
    "V" = 1.85
    "W" = 0.6000000000000001
    "X" = -1
    "Y" = -1
    "Z" = 0.8
    "AA" = -1
    "AB" = -0.5
    "AC" = 0.425
}
Set-RowData 182 $d182

# --- Row 183: brand-new match row ---
$d183 = [ordered]@{
    "A" = 181
    "B" = 8050912
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Apertura"
    "E" = 45398.75
    "F" = "Montevideo Wanderers"
    "G" = "Liverpool Montevideo"
    "H" = 2
    "I" = 3
    "J" = "A"
    "K" = 3.2
    "L" = 3.3
    "M" = 2.2
    "N" = 3.25
    "O" = 3.4
    "P" = 2.15
    "Q" = 0.25
    "R" = 1.975
    "S" = 1.875
    "T" = 2.25
    "U" = 1.825
    "V" = 2.025
    "W" = -1
    "X" = -1
    "Y" = 1.15
    "Z" = -1
    "AA" = 0.875
    "AB" = 0.825
    "AC" = -1
}
Set-RowData 183 $d183
# Copy the number-format / border styling used on the "id" (col A) and "Date" (col E) columns
$ws.Range("A114").Copy()
$ws.Range("A183").PasteSpecial(-4122)
$ws.Range("E114").Copy()
$ws.Range("E183").PasteSpecial(-4122)

# --- Row 184: brand-new match row ---
$d184 = [ordered]@{
    "A" = 182
    "B" = 8081163
    "C" = "Uruguay Primera División"
    "D" = "Uruguay Apertura"
    "E" = 45401.5625
    "F" = "Racing Club de Montevideo"
    "G" = "Cerro"
    "K" = 2.3
    "L" = 3.2
    "M" = 3.2
    "N" = 2.25
    "O" = 3.2
    "P" = 3.25
    "Q" = -0.25
    "R" = 1.975
    "S" = 1.875
    "T" = 2.25
    "U" = 2.025
    "V" = 1.825
    "W" = 0
    "X" = 0
    "Y" = 0
    "Z" = 0
    "AA" = 0
}
Set-RowData 184 $d184
# Copy the number-format / border styling used on the "id" (col A) and "Date" (col E) columns
$ws.Range("A114").Copy()
$ws.Range("A184").PasteSpecial(-4122)
$ws.Range("E114").Copy()
$ws.Range("E184").PasteSpecial(-4122)

